# Fill in vocabulary list (column A) and mark which column contains the
# "15/04/2021" date string in each of the first six rows, mirroring the
# "loop to read words" exercise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A values, entered in the order that reproduces the original
#     shared-string table ordering (1,2,3,6,5,4,date,7..16) ---
$ws.Range("A1").Value = "sample 1"
$ws.Range("A2").Value = "sample 2"
$ws.Range("A3").Value = "sample 3"
$ws.Range("A6").Value = "sample 6"
$ws.Range("A5").Value = "sample 5"
$ws.Range("A4").Value = "sample 4"

# --- Create the date-format style once on G1 (this also introduces the
#     "15/04/2021" shared string) ---
$ws.Range("G1").NumberFormat = "mm-dd-yy"
$ws.Range("G1").Value = "15/04/2021"
$ws.Range("G1").Copy()

# --- Reuse that exact style (format-only paste) on every other cell in the
#     grid that needs it: the "blank" marker cells plus the other cells
#     that will hold the date text ---
$ws.Range("B1:F1").PasteSpecial(-4122)

$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D2:G2").PasteSpecial(-4122)

$ws.Range("B3:D3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("F3:G3").PasteSpecial(-4122)

$ws.Range("B4:D4").PasteSpecial(-4122)
$ws.Range("F4:G4").PasteSpecial(-4122)

$ws.Range("B5:G5").PasteSpecial(-4122)

$ws.Range("B6:G6").PasteSpecial(-4122)

$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)

# --- Now write the date text into the remaining cells (style is already
#     applied, so this does not create extra cellXfs entries) ---
$ws.Range("F2").Value = "15/04/2021"
$ws.Range("E3").Value = "15/04/2021"
$ws.Range("D4").Value = "15/04/2021"
$ws.Range("C5").Value = "15/04/2021"
$ws.Range("C14").Value = "15/04/2021"
$ws.Range("B15").Value = "15/04/2021"

# --- Remaining word list, rows 7-16 ---
$ws.Range("A7").Value = "sample 7"
$ws.Range("A8").Value = "sample 8"
$ws.Range("A9").Value = "sample 9"
$ws.Range("A10").Value = "sample 10"
$ws.Range("A11").Value = "sample 11"
$ws.Range("A12").Value = "sample 12"
$ws.Range("A13").Value = "sample 13"
$ws.Range("A14").Value = "sample 14"
$ws.Range("A15").Value = "sample 15"
$ws.Range("A16").Value = "sample 16"

# --- Column widths (best-fit look of the original workbook) ---
$ws.Columns("A").ColumnWidth = 9
$ws.Columns("B").ColumnWidth = 8.8
$ws.Columns("C:G").ColumnWidth = 9.75

# --- Final selected cell, matching the author's last cursor position ---
[void]$ws.Range("L18").Select()
